$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALERT")

# New row 5 values (mirrors the ALERT_002/ALERT_003 rows above it)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "ALERT_004"
$ws.Range("C5").Value = "C:\\Git_Evergreen\\fms_cba\\DataSet\\Integration_DataSet\\Extracts\\DNR\\DNR_Reports\\"
$ws.Range("D5").Value = "|"
$ws.Range("E5").Value = "Outstanding_4"
$ws.Range("F5").Value = "Deal Name|Deal Tracking Number|Alias Number|Alert Heading|Alert Content|User Name|Date Added / Amended"

# Match formatting used by the rest of the table (Text number format, Arial 10)
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Font.Name = "Arial"
$ws.Range("C5").Font.Size = 10

$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Font.Name = "Arial"
$ws.Range("F5").Font.Size = 10

# Update the active selection to reflect the new last cell, as in the edited file
$ws.Range("F6").Select()
